$d = $word.ActiveDocument

# Small helper fragments used with Range.InsertXML to splice raw OOXML
# markup (so we can produce a genuinely-empty <w:p/> or drop in elements,
# like <w:lastRenderedPageBreak/>, that have no COM property of their own).
$pkgOpen = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Number of paragraphs in the document before this edit.
$orig = $d.Paragraphs.Count

# Append 7 new paragraphs after the current last paragraph (all of them
# land after it, and before the section properties):
#   orig+1, orig+2, orig+3  -> empty
#   orig+4                  -> "Environment" heading
#   orig+5                  -> empty
#   orig+6                  -> "Understanding Unity-ML Agent's Environment:"
#   orig+7                  -> "This project uses "
$p = $d.Paragraphs($orig).Range
for ($i = 0; $i -lt 7; $i++) {
    $p.InsertParagraphAfter()
    $p = $d.Paragraphs($orig + $i + 1).Range
}

# orig+1 .. orig+3: make them truly empty paragraphs (no stray empty run).
for ($i = 1; $i -le 3; $i++) {
    $blankRange = $d.Paragraphs($orig + $i).Range
    $blankRange.InsertXML($pkgOpen + '<w:p/>' + $pkgClose) | Out-Null
}

# orig+4: "Environment" heading. Word stamps a lastRenderedPageBreak on
# the run because it starts the next printed page.
$envParaRange = $d.Paragraphs($orig + 4).Range
$envParaRange.InsertXML($pkgOpen + '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Environment</w:t></w:r></w:p>' + $pkgClose) | Out-Null

# orig+5: empty paragraph again.
$blankRange = $d.Paragraphs($orig + 5).Range
$blankRange.InsertXML($pkgOpen + '<w:p/>' + $pkgClose) | Out-Null

# orig+6: "Understanding Unity-ML Agent's Environment:"
$understandPara = $d.Paragraphs($orig + 6).Range
$understandPara.Text = "Understanding Unity-ML Agent" + [char]0x2019 + "s Environment:"

# orig+7 (now the document's final paragraph): "This project uses "
$lastPara = $d.Paragraphs($orig + 7).Range
$lastPara.Text = "This project uses "
